# Update the results sheet: rewrite row 2 (Steve Coup), row 3 (Fredrick Ndote, scores
# only) and row 4 (Bostwald Kite, replacing the old "Steve Coup" row), then remove the
# two trailing rows (old row 5 "Bostwald Kite" and old row 6 "Nderi Neti") that are no
# longer part of the dataframe that was written out.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Delete the last two data rows first (rows 5 and 6) so the remaining rows keep their
# original row numbers while we overwrite them.
$ws.Rows(6).EntireRow.Delete() | Out-Null
$ws.Rows(5).EntireRow.Delete() | Out-Null

# Row 2: now holds "Steve Coup" (previously file #15615 / Heri Dati)
$ws.Cells.Item(2, 1).Value = 15613
$ws.Cells.Item(2, 2).Value = "Steve Coup"
$ws.Cells.Item(2, 3).Value = 80
$ws.Cells.Item(2, 4).Value = 44
$ws.Cells.Item(2, 5).Value = 68
$ws.Cells.Item(2, 6).Value = 48
$ws.Cells.Item(2, 7).Value = 70
$ws.Cells.Item(2, 8).Value = 81
$ws.Cells.Item(2, 9).Value = 54
$ws.Cells.Item(2, 10).Value = 89
$ws.Cells.Item(2, 11).Value = 534
$ws.Cells.Item(2, 12).Value = 66.75
$ws.Cells.Item(2, 13).Value = "B-"
$ws.Cells.Item(2, 14).Value = 1

# Row 3: "Fredrick Ndote" keeps his file number/name, only scores change
$ws.Cells.Item(3, 3).Value = 70
$ws.Cells.Item(3, 4).Value = 64
$ws.Cells.Item(3, 5).Value = 50
$ws.Cells.Item(3, 6).Value = 59
$ws.Cells.Item(3, 7).Value = 90
$ws.Cells.Item(3, 8).Value = 76
$ws.Cells.Item(3, 9).Value = 36
$ws.Cells.Item(3, 10).Value = 86
$ws.Cells.Item(3, 11).Value = 531
$ws.Cells.Item(3, 12).Value = 66.375
$ws.Cells.Item(3, 13).Value = "B-"

# Row 4: now holds "Bostwald Kite" (previously file #15613 / Steve Coup)
$ws.Cells.Item(4, 1).Value = 15612
$ws.Cells.Item(4, 2).Value = "Bostwald Kite"
$ws.Cells.Item(4, 3).Value = 56
$ws.Cells.Item(4, 4).Value = 79
$ws.Cells.Item(4, 5).Value = 36
$ws.Cells.Item(4, 6).Value = 60
$ws.Cells.Item(4, 7).Value = 87
$ws.Cells.Item(4, 8).Value = 48
$ws.Cells.Item(4, 9).Value = 68
$ws.Cells.Item(4, 10).Value = 69
$ws.Cells.Item(4, 11).Value = 503
$ws.Cells.Item(4, 12).Value = 62.875
$ws.Cells.Item(4, 13).Value = "C+"
$ws.Cells.Item(4, 14).Value = 3
